# Helper: convert a "RRGGBB" hex color string into the long integer value
# used by the PowerPoint COM RGB()-style color properties (R + G*256 + B*65536).
function ToVbaRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5's table (graphic frame "Google Shape;122;p17") switches to a
#    different built-in table style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{35EA6AD4-2CF2-463D-814F-88C7CBC9745E}")

# ---------------------------------------------------------------------------
# 2) The deck's main theme (slide master -> theme1.xml, currently the
#    "Integral" / Red Violet palette) is recolored to the default Office
#    palette.
# ---------------------------------------------------------------------------
$officePalette = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$mainScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $mainScheme.Colors($i).RGB = ToVbaRgb($officePalette[$i - 1])
}

# ---------------------------------------------------------------------------
# 3) The Notes Master's theme (theme2.xml, currently the default Office
#    palette) is recolored to the "Integral" / Red Violet palette that used
#    to live on the main deck theme.
# ---------------------------------------------------------------------------
$redVioletPalette = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "454551",  # dk2
    "D8D9DC",  # lt2
    "E32D91",  # accent1
    "C830CC",  # accent2
    "4EA6DC",  # accent3
    "4775E7",  # accent4
    "8971E1",  # accent5
    "D54773",  # accent6
    "6B9F25",  # hlink
    "8C8C8C"   # folHlink
)

$notesScheme = $p.NotesMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Colors($i).RGB = ToVbaRgb($redVioletPalette[$i - 1])
}
